$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.703.87"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "3.646.25"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  +16.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "656.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.421"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").Value = "3.644.70"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "4.327.28"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "96.546.22"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000258"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.628.31"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.526"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.88%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "513.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("E29").Value = "  +11.46%  "
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "614.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "43.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +26.35%  "
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("E42").Value = "  +5.85%  "
$ws.Range("E43").Value = "  +6.42%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0439"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.403"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.49%  "
